$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Make sure "ODI Batting" is not the active sheet while we copy FROM it further
# down (copying from the active sheet behaves unreliably in this engine).
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()

# 1. Clear the (already empty) B3 / B5 / B6 cells on "ODI Batting" so they
#    disappear from the sheet entirely, matching the target workbook.
$odiBatting.Range("B3").ClearContents()
$odiBatting.Range("B5").ClearContents()
$odiBatting.Range("B6").ClearContents()

# 2. Add the new "ODI Batting Extra" worksheet as the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"
$newSheet.Activate()

# Match the page margins used by the rest of the workbook's sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy the header row formatting/style from "ODI Batting" (re-using the same
# header style as the other sheets) and then overwrite the header text.
$odiBatting.Range("A1:F1").Copy($newSheet.Range("A1:F1"))

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Copy the MATCH_CODE values (already stored as text) from the "ODI Batting"
# sheet so they keep their text data type.
$odiBatting.Range("D2:D6").Copy($newSheet.Range("A2:A6"))

# BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL are left blank except
# for two numeric BATTING_POSITION values.
$newSheet.Range("B5").Value = 10
$newSheet.Range("B6").Value = 10

# MAN_OF_MATCH column.
$newSheet.Range("F2").Value = "NO"
$newSheet.Range("F3").Value = "NO"
$newSheet.Range("F4").Value = "NO"
$newSheet.Range("F5").Value = "NO"
$newSheet.Range("F6").Value = "NO"
